$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the "RUNDOWN TUGAS :" bullet list by its known, stable
# starting paragraph text, then work with indices relative to it so the
# script is not order-of-edit fragile.
# ---------------------------------------------------------------------
$base = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Contains("RUNDOWN TUGAS")) {
        $base = $i
        break
    }
}

# Bullets, relative to $base (1-indexed offsets):
#  base+1  Rubah sistem tarif flat menjadi 3 tahap perincian dan perubahan biaya registrasi. Rubah harga bisa dilakukan di dalam konfigurasi agar suatu saat mudah terjadi perubahan harga.
#  base+2  Rubah denda juga bisa dikonfigurasi di sistem.
#  base+3  Rubah sistem jatuh tempo. Sehingga bisa diedit di konfigurasi.
#  base+4  Rubah pengetikan dusun tidak manual, tapi menggunakan combo box. ...
#  base+5  ID transaksi di pembelian barang dibuat otomatis, view id dihilangkan saja.
#  base+6  ID pengeluaran juga dibuat otomatis. ...
#  base+7  ID Barang juga dimasukan ke dalam sistem.
#  base+8  Setelah 3 bulan denda, diputus.
#  base+9  Buat form informasi keluhan untuk memudahkan petugas mengatur jadwal perbaikan saluran air.
#  base+10 Perbaiki tampilan home biar ga kosong.   <- has the _GoBack bookmark at its end
#  base+11 Neraca.

# ---------------------------------------------------------------------
# 1. Move the _GoBack bookmark from the end of "Perbaiki tampilan
#    home..." (base+10) to the start of the bullet that will end up
#    reading "Rubah pengetikan dusun..." (that is base+3, once base+2
#    is deleted below) *before* any text is shifted around.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRng = $d.Paragraphs($base + 3).Range.Duplicate
$bmRng.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

# ---------------------------------------------------------------------
# 2. Delete the whole "Perbaiki tampilan home..." bullet (base+10); its
#    text re-appears, shifted, further up the list below.
# ---------------------------------------------------------------------
$d.Paragraphs($base + 10).Range.Delete()

# ---------------------------------------------------------------------
# 3. Shift bullet texts: base+9 down to base+4 each take over what used
#    to be the text of the *next* bullet (base+10 .. base+5). Go from
#    the bottom up so we always read an not-yet-modified paragraph.
# ---------------------------------------------------------------------
$d.Paragraphs($base + 9).Range.Text  = "Perbaiki tampilan home biar ga kosong."
$d.Paragraphs($base + 8).Range.Text  = "Buat form informasi keluhan untuk memudahkan petugas mengatur jadwal perbaikan saluran air."
$d.Paragraphs($base + 7).Range.Text  = "Setelah 3 bulan denda, diputus."
$d.Paragraphs($base + 6).Range.Text  = "ID Barang juga dimasukan ke dalam sistem."
$d.Paragraphs($base + 5).Range.Text  = "ID pengeluaran juga dibuat otomatis. Tampilan dihilangkan proses pengetikan ID. (sama sistemnya dgn no 3)"
$d.Paragraphs($base + 4).Range.Text  = "ID transaksi di pembelian barang dibuat otomatis, view id dihilangkan saja."

# ---------------------------------------------------------------------
# 4. The bullet at base+3 ("Rubah sistem jatuh tempo...") now takes the
#    text that used to belong to base+4 ("Rubah pengetikan dusun..."),
#    continuing the shift pattern; it is also the paragraph that just
#    received the relocated _GoBack bookmark above.
# ---------------------------------------------------------------------
$d.Paragraphs($base + 3).Range.Text = "Rubah pengetikan dusun tidak manual, tapi menggunakan combo box. Kecamatan dll ditampilkan secara otomatis untuk memudahkan registrasi data."

# ---------------------------------------------------------------------
# 5. Delete bullet base+2 ("Rubah denda juga bisa dikonfigurasi di
#    sistem.") outright - its own content never survives the shift.
# ---------------------------------------------------------------------
$d.Paragraphs($base + 2).Range.Delete()

# ---------------------------------------------------------------------
# 6. Shorten bullet base+1: drop the second sentence, keep trailing
#    space.
# ---------------------------------------------------------------------
$d.Paragraphs($base + 1).Range.Text = "Rubah sistem tarif flat menjadi 3 tahap perincian dan perubahan biaya registrasi. "

# ---------------------------------------------------------------------
# 7. Insert a brand new bullet before base+1 with the new text, in red.
# ---------------------------------------------------------------------
$insertRng = $d.Paragraphs($base + 1).Range.Duplicate
$insertRng.Collapse(1)
$insertRng.InsertBefore("Rubah semua tarif (abodemen,pertama,kedua,ketiga,jatuh tempo, denda) bisa diseting untuk jangka panjang dan pengecekan format angka.`r")
$d.Paragraphs($base + 1).Range.Font.Color = 255
